# Update "想去人数" (want-to-go count) figures in the F column across the
# four sheets of the workbook to match the freshly generated site output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 26331
$ws1.Range("F4").Value  = 577
$ws1.Range("F5").Value  = 250
$ws1.Range("F6").Value  = 592
$ws1.Range("F8").Value  = 531
$ws1.Range("F10").Value = 349
$ws1.Range("F11").Value = 219
$ws1.Range("F12").Value = 181
$ws1.Range("F13").Value = 46
$ws1.Range("F14").Value = 293
$ws1.Range("F15").Value = 47
$ws1.Range("F16").Value = 369
$ws1.Range("F18").Value = 1506
$ws1.Range("F19").Value = 183
$ws1.Range("F20").Value = 29
$ws1.Range("F21").Value = 428
$ws1.Range("F22").Value = 99
$ws1.Range("F23").Value = 115

# --- Sheet "演出" (rId2 / sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 227
$ws2.Range("F6").Value  = 181
$ws2.Range("F7").Value  = 36
$ws2.Range("F8").Value  = 109
$ws2.Range("F9").Value  = 109
$ws2.Range("F10").Value = 433
$ws2.Range("F15").Value = 55
$ws2.Range("F16").Value = 24
$ws2.Range("F17").Value = 27

# --- Sheet "本地生活" (rId3 / sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 4991
$ws3.Range("F3").Value = 207

# --- Sheet "全部类型" (rId4 / sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 4991
$ws4.Range("F4").Value  = 207
$ws4.Range("F5").Value  = 26331
$ws4.Range("F6").Value  = 577
$ws4.Range("F8").Value  = 250
$ws4.Range("F9").Value  = 227
$ws4.Range("F10").Value = 592
$ws4.Range("F14").Value = 181
$ws4.Range("F15").Value = 181
$ws4.Range("F16").Value = 36
$ws4.Range("F17").Value = 109
$ws4.Range("F18").Value = 109
$ws4.Range("F19").Value = 433
$ws4.Range("F20").Value = 531
$ws4.Range("F23").Value = 349
$ws4.Range("F24").Value = 219
$ws4.Range("F25").Value = 181
$ws4.Range("F26").Value = 46
$ws4.Range("F28").Value = 293
$ws4.Range("F29").Value = 47
$ws4.Range("F32").Value = 369
$ws4.Range("F34").Value = 55
$ws4.Range("F35").Value = 1506
$ws4.Range("F36").Value = 183
$ws4.Range("F37").Value = 24
$ws4.Range("F38").Value = 29
$ws4.Range("F39").Value = 428
$ws4.Range("F40").Value = 99
$ws4.Range("F41").Value = 115
$ws4.Range("F42").Value = 27
